$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-08-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-03 Thursday", 2) | Out-Null

# Update each multiplication-table cell by position (row, column) to avoid
# any ambiguity from duplicate/overlapping values produced during the edit.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "49×32="
$t.Cell(1, 2).Range.Text = "97×16="
$t.Cell(1, 3).Range.Text = "84×66="
$t.Cell(1, 4).Range.Text = "100×53="
$t.Cell(1, 5).Range.Text = "93×56="
$t.Cell(2, 1).Range.Text = "93×92="
$t.Cell(2, 2).Range.Text = "67×82="
$t.Cell(2, 3).Range.Text = "12×38="
$t.Cell(2, 4).Range.Text = "56×34="
$t.Cell(2, 5).Range.Text = "15×34="
$t.Cell(3, 1).Range.Text = "70×62="
$t.Cell(3, 2).Range.Text = "98×35="
$t.Cell(3, 3).Range.Text = "59×85="
$t.Cell(3, 4).Range.Text = "35×29="
$t.Cell(3, 5).Range.Text = "93×46="
$t.Cell(4, 1).Range.Text = "57×11="
$t.Cell(4, 2).Range.Text = "14×68="
$t.Cell(4, 3).Range.Text = "17×84="
$t.Cell(4, 4).Range.Text = "70×15="
$t.Cell(4, 5).Range.Text = "86×73="
$t.Cell(5, 1).Range.Text = "10×79="
$t.Cell(5, 2).Range.Text = "72×69="
$t.Cell(5, 3).Range.Text = "94×62="
$t.Cell(5, 4).Range.Text = "60×87="
$t.Cell(5, 5).Range.Text = "76×33="
$t.Cell(6, 1).Range.Text = "78×28="
$t.Cell(6, 2).Range.Text = "88×89="
$t.Cell(6, 3).Range.Text = "54×66="
$t.Cell(6, 4).Range.Text = "17×66="
$t.Cell(6, 5).Range.Text = "46×63="
$t.Cell(7, 1).Range.Text = "87×59="
$t.Cell(7, 2).Range.Text = "99×93="
$t.Cell(7, 3).Range.Text = "84×81="
$t.Cell(7, 4).Range.Text = "15×99="
$t.Cell(7, 5).Range.Text = "13×34="
$t.Cell(8, 1).Range.Text = "88×99="
$t.Cell(8, 2).Range.Text = "39×16="
$t.Cell(8, 3).Range.Text = "52×80="
$t.Cell(8, 4).Range.Text = "86×36="
$t.Cell(8, 5).Range.Text = "55×55="
$t.Cell(9, 1).Range.Text = "94×69="
$t.Cell(9, 2).Range.Text = "41×13="
$t.Cell(9, 3).Range.Text = "97×70="
$t.Cell(9, 4).Range.Text = "41×63="
$t.Cell(9, 5).Range.Text = "81×73="
$t.Cell(10, 1).Range.Text = "78×95="
$t.Cell(10, 2).Range.Text = "83×41="
$t.Cell(10, 3).Range.Text = "34×36="
$t.Cell(10, 4).Range.Text = "37×52="
$t.Cell(10, 5).Range.Text = "78×29="
$t.Cell(11, 1).Range.Text = "62×51="
$t.Cell(11, 2).Range.Text = "46×60="
$t.Cell(11, 3).Range.Text = "51×59="
$t.Cell(11, 4).Range.Text = "20×80="
$t.Cell(11, 5).Range.Text = "67×26="
$t.Cell(12, 1).Range.Text = "66×61="
$t.Cell(12, 2).Range.Text = "98×92="
$t.Cell(12, 3).Range.Text = "12×16="
$t.Cell(12, 4).Range.Text = "73×81="
$t.Cell(12, 5).Range.Text = "97×38="
$t.Cell(13, 1).Range.Text = "18×55="
$t.Cell(13, 2).Range.Text = "100×76="
$t.Cell(13, 3).Range.Text = "40×100="
$t.Cell(13, 4).Range.Text = "96×40="
$t.Cell(13, 5).Range.Text = "95×46="
$t.Cell(14, 1).Range.Text = "74×86="
$t.Cell(14, 2).Range.Text = "10×24="
$t.Cell(14, 3).Range.Text = "90×72="
$t.Cell(14, 4).Range.Text = "28×64="
$t.Cell(14, 5).Range.Text = "88×63="
$t.Cell(15, 1).Range.Text = "58×49="
$t.Cell(15, 2).Range.Text = "54×69="
$t.Cell(15, 3).Range.Text = "79×56="
$t.Cell(15, 4).Range.Text = "91×26="
$t.Cell(15, 5).Range.Text = "54×72="
$t.Cell(16, 1).Range.Text = "59×21="
$t.Cell(16, 2).Range.Text = "17×14="
$t.Cell(16, 3).Range.Text = "54×43="
$t.Cell(16, 4).Range.Text = "83×19="
$t.Cell(16, 5).Range.Text = "82×42="
$t.Cell(17, 1).Range.Text = "72×31="
$t.Cell(17, 2).Range.Text = "54×64="
$t.Cell(17, 3).Range.Text = "28×87="
$t.Cell(17, 4).Range.Text = "74×37="
$t.Cell(17, 5).Range.Text = "25×66="
$t.Cell(18, 1).Range.Text = "30×92="
$t.Cell(18, 2).Range.Text = "78×76="
$t.Cell(18, 3).Range.Text = "77×88="
$t.Cell(18, 4).Range.Text = "72×75="
$t.Cell(18, 5).Range.Text = "96×52="
$t.Cell(19, 1).Range.Text = "69×27="
$t.Cell(19, 2).Range.Text = "72×51="
$t.Cell(19, 3).Range.Text = "82×40="
$t.Cell(19, 4).Range.Text = "56×21="
$t.Cell(19, 5).Range.Text = "74×23="
$t.Cell(20, 1).Range.Text = "40×84="
$t.Cell(20, 2).Range.Text = "61×20="
$t.Cell(20, 3).Range.Text = "42×19="
$t.Cell(20, 4).Range.Text = "85×69="
$t.Cell(20, 5).Range.Text = "87×28="
